$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder row 2 values for columns A-F (G:J remain unchanged)
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2

# Update the active selection to J2
$ws.Range("J2").Select()
